$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The three target cells (C2, C3, C24) hold numeric-looking Service IDs that
# must be stored as TEXT (shared string), matching how the rest of the sheet
# stores this kind of value. Assigning a numeric-looking string straight to
# .Value causes Excel to auto-convert it to a real number, so instead we
# stage the text in an out-of-the-way helper cell that is explicitly
# formatted as Text, copy it, and paste-special (values only) into the
# target cell. That preserves the target cell's original (default) style
# while still writing the value as text.

$helper = $ws.Range("Z500")
$helper.NumberFormat = "@"

$helper.Value = "10304774"
$helper.Copy()
$ws.Range("C2").PasteSpecial(-4163)   # xlPasteValues

$helper.Value = "10304776"
$helper.Copy()
$ws.Range("C3").PasteSpecial(-4163)   # xlPasteValues

$helper.Value = "137626708"
$helper.Copy()
$ws.Range("C24").PasteSpecial(-4163)  # xlPasteValues

# Remove the helper cell/row remnants so the sheet dimensions/content stay
# exactly as before outside of the three intended cells.
$helper.Delete(-4162)  # xlShiftUp

$excel.CutCopyMode = $false
